# Append the new Oct/Nov-2025 "ORA Errors" data point as row 24:
#   A24 = 12/1/2025 (date, formatted like the rest of column A)
#   B24 = 76        (plain error-count number, like the rest of column B)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format from the last existing date cell (A23) down onto A24 first,
# so the new cell inherits the same date number-format style (style index 1)
# instead of Excel fabricating a brand-new custom number format.
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)   # xlPasteFormats

# Now write the actual values for the new row.
$ws.Range("A24").Value = 45992          # serial date -> 12/1/2025
$ws.Range("B24").Value = 76

# Match the saved selection state: A24:B24 selected, A24 active.
$ws.Range("A24:B24").Select()
